$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column R for year 2021 by copying column Q (2020) and inserting
# the copy to the right, so it inherits the same cell formatting/styles.
$ws.Range("Q1:Q13").Copy() | Out-Null
$ws.Range("R1:R13").Insert(-4161) | Out-Null  # xlShiftToRight

# Overwrite the copied values with the actual 2021 figures.
$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 13.5
$ws.Range("R5").Value = 15.1

# Update the active selection to match the edited workbook.
$ws.Range("T3").Select() | Out-Null
